$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddress, $value) {
    $r = $ws.Range($cellAddress)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "29.866.59"
$ws.Range("E2").Value = "  +2.74%  "

Set-TextValue "D3" "1.860.17"
$ws.Range("E3").Value = "  +2.25%  "

Set-TextValue "D4" "1.000"
$ws.Range("E4").Value = "  +0.00%  "

Set-TextValue "D6" "0.6363"
$ws.Range("E6").Value = "  +3.53%  "

$ws.Range("E7").Value = "  +0.03%  "

Set-TextValue "D8" "0.2999"
$ws.Range("E8").Value = "  +4.15%  "

Set-TextValue "D9" "0.07476"
$ws.Range("E9").Value = "  +2.12%  "

Set-TextValue "D10" "24.49"
$ws.Range("E10").Value = "  +7.21%  "

Set-TextValue "D11" "0.07682"
$ws.Range("E11").Value = "  +0.36%  "

Set-TextValue "D12" "1.856.41"
$ws.Range("E12").Value = "  +2.03%  "

Set-TextValue "D13" "5.052"
$ws.Range("E13").Value = "  +2.30%  "

Set-TextValue "D14" "0.6894"
$ws.Range("E14").Value = "  +4.65%  "

Set-TextValue "D15" "84.35"
$ws.Range("E15").Value = "  +3.41%  "

Set-TextValue "D16" "0.000009396"
$ws.Range("E16").Value = "  +4.11%  "

Set-TextValue "D17" "6.095"
$ws.Range("E17").Value = "  +4.68%  "

Set-TextValue "D18" "29.837.60"
$ws.Range("E18").Value = "  +2.75%  "

Set-TextValue "D19" "2.116.83"
$ws.Range("E19").Value = "  +2.45%  "

Set-TextValue "D20" "238.86"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("E21").Value = "  +2.05%  "

$ws.Range("E22").Value = "  -0.01%  "

Set-TextValue "D23" "7.349"
$ws.Range("E23").Value = "  +3.42%  "

$ws.Range("E24").Value = "  +0.09%  "

Set-TextValue "D25" "159.23"
$ws.Range("E25").Value = "  +1.02%  "

$ws.Range("E26").Value = "  +1.19%  "

Set-TextValue "D27" "8.580"
$ws.Range("E27").Value = "  +1.93%  "

Set-TextValue "D28" "17.97"
$ws.Range("E28").Value = "  +2.31%  "

Set-TextValue "D29" "1.504"
$ws.Range("E29").Value = "  +1.35%  "

Set-TextValue "D30" "0.06064"
$ws.Range("E30").Value = "  +9.07%  "

Set-TextValue "D31" "1.273"
$ws.Range("E31").Value = "  +5.39%  "

$ws.Range("E32").Value = "  +1.33%  "

$ws.Range("E33").Value = "  +1.39%  "

Set-TextValue "D34" "1.895"
$ws.Range("E34").Value = "  +4.75%  "

Set-TextValue "D35" "1.166"
$ws.Range("E35").Value = "  +3.24%  "

Set-TextValue "D36" "0.7283"
$ws.Range("E36").Value = "  -0.62%  "

Set-TextValue "D37" "2.611"
$ws.Range("E37").Value = "  -0.04%  "

Set-TextValue "D38" "2.855"
$ws.Range("E38").Value = "  +0.87%  "

Set-TextValue "D39" "0.01795"
$ws.Range("E39").Value = "  +2.70%  "

Set-TextValue "D40" "1.222.35"
$ws.Range("E40").Value = "  +1.27%  "

$ws.Range("E41").Value = "  +4.29%  "

Set-TextValue "D42" "6.303"
$ws.Range("E42").Value = "  -0.67%  "

Set-TextValue "D43" "1.002"
$ws.Range("E43").Value = "  +0.09%  "

Set-TextValue "D44" "2.021.44"
$ws.Range("E44").Value = "  +2.64%  "

Set-TextValue "D45" "102.35"
$ws.Range("E45").Value = "  +1.51%  "

Set-TextValue "D46" "66.31"
$ws.Range("E46").Value = "  +3.03%  "

Set-TextValue "D47" "0.00000000122"
$ws.Range("E47").Value = "  +3.80%  "

Set-TextValue "D48" "0.5092"
$ws.Range("E48").Value = "  +0.17%  "

Set-TextValue "D49" "9.274"
$ws.Range("E49").Value = "  +2.63%  "

Set-TextValue "D50" "0.4089"
$ws.Range("E50").Value = "  +2.49%  "

Set-TextValue "D51" "0.1143"
$ws.Range("E51").Value = "  +3.28%  "
